# Regenerate sval data to filter save games: update computed columns B:G
# for rows 2-5 on the active worksheet with the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.7287194209349384, 0.3375848360084654, 3.082599426703578, 6.48142807727062, 1, 10.6303317609176)
    3 = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 4.371470058157054)
    4 = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 0, 10.35301142835362)
    5 = @(0.7287194209349384, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 3.594575437922795)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $col = 2
    foreach ($v in $vals) {
        $ws.Cells.Item($row, $col).Value = $v
        $col = $col + 1
    }
}
